# Daily attendance processing - 2025-12-06 20:50:27
#
# Normalises the "Recorded By" column (G): whenever the recorded-by list is a
# two-name pair ending in "System" (e.g. "dnasr281@gmail.com, System"), move
# "System" to the front (e.g. "System, dnasr281@gmail.com"). Rows already
# starting with "System", rows without "System" at all, and rows listing more
# than two names are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

$col = 7  # column G - "Recorded By"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($val -ne $null) {
        $parts = $val -split ", "
        if ($parts.Count -eq 2 -and $parts[1] -eq "System") {
            $cell.Value = "System, " + $parts[0]
        }
    }
}
